$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.447.09'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.569.72'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.08'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('E10').Value = '  -0.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0866'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = '1.795.34'
$ws.Range('D13').Value = '1.565.63'
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.520'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').Value = '27.442.64'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '213.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.72%  '
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.25'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.70'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('E24').Value = '  +1.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.03'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('E32').Value = '  -1.91%  '
$ws.Range('D33').Value = '1.358.13'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.822'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.25%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.972'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.18%  '
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.14'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.30'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = '1.708.67'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.32'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('D49').Value = '0.0₇0997'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0954'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.64%  '
$ws.Range('E51').Value = '  -0.66%  '
